$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.167.89"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.639.30"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.991"
$ws.Range("E4").Value = "  -1.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.70"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.993"
$ws.Range("E7").Value = "  -0.88%  "
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0634"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "1.866.80"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "1.619.30"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.35"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "26.136.40"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.993"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.46"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.87"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.00"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.37"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.993"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.78"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.90"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.66"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.59"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.907"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").Value = "1.144.02"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.547"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.992"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.19"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.792"
$ws.Range("E44").Value = "  -2.02%  "
$ws.Range("D45").Value = "1.776.69"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0107"
$ws.Range("E46").Value = "  -4.30%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.86"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("E49").Value = "  +4.99%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("E51").Value = "  +3.11%  "
